# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates DAMSLTag (col I) and DialogAct (col J)
# for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;  DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 13; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 29; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 30; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 47; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 64; DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 67; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 79; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 89; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.DAMSLTag
    $ws.Range("J" + $u.Row).Value = $u.DialogAct
}
